$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.405.39"
$ws.Range("D2").Style = $ws.Range("B2").Style
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.56"
$ws.Range("D3").Style = $ws.Range("B3").Style
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.66"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6324"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07568"
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.58"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.58"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  -5.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.005"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6805"
$ws.Range("D14").Style = $ws.Range("B14").Style
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("E15").Value = "  +5.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.42"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.110.97"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "  -6.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.175"
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.438.00"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.04"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  -1.27%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.492"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.79"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1395"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.353"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.61"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.461"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.301"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +3.45%  "
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.104"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.026"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.851"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = "  -0.28%  "
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7106"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.594"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.250.88"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.771"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  -1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.382"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9015"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.80"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.94"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000119"
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.100"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4004"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.675"
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.928"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("E51").Value = "  -0.08%  "
